$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 updates
$ws.Range("Q2").Value = 2.3
$ws.Range("R2").Value = 1.6
$ws.Range("BD2").Value = 126

# Row 4 updates
$ws.Range("I4").Value = 4.5
$ws.Range("L4").Value = 5
$ws.Range("M4").Value = 1.07
$ws.Range("N4").Value = 9
$ws.Range("O4").Value = 1.36
$ws.Range("P4").Value = 3
$ws.Range("Q4").Value = 2.2
$ws.Range("R4").Value = 1.65
$ws.Range("Z4").Value = 15
$ws.Range("AG4").Value = 1000
$ws.Range("AI4").Value = 21
$ws.Range("AJ4").Value = 15
$ws.Range("AW4").Value = 6
$ws.Range("AX4").Value = 26
